$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report row was inserted above the old row 100, pushing
# the existing rows 100-107 down to 101-108.
$ws.Rows(100).Insert()

$ws.Range("A100").Value = 7
$ws.Range("B100").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C100").Value = "Ñuble"
$ws.Range("D100").Value = 44918
$ws.Range("E100").Value = 16
$ws.Range("F100").Value = 100112031
$ws.Range("G100").Value = "Poroto verde"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 60
$ws.Range("K100").Value = 27000
$ws.Range("L100").Value = 28000
$ws.Range("M100").Value = 27500
$ws.Range("N100").Value = "$/saco 25 kilos"
$ws.Range("O100").Value = "Región del Maule"
$ws.Range("P100").Value = 1100
$ws.Range("Q100").Value = 25
$ws.Range("R100").Value = "Hortaliza"
